$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff. D-column values that
# look like plain numbers ("414.40", "11.00", ...) must stay stored as text
# (matching the original inlineStr cells), so we force NumberFormat = "@"
# (Text) on just those specific cells before writing the value - this stops
# Excel from auto-converting them to numeric values.

$ws.Range("D2").Value = "62.635.79"
$ws.Range("E2").Value = "  +5.93%  "

$ws.Range("D3").Value = "3.464.97"
$ws.Range("E3").Value = "  +3.65%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.40"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.60"
$ws.Range("E6").Value = "  +16.70%  "

$ws.Range("D7").Value = "3.456.19"
$ws.Range("E7").Value = "  +3.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  +1.36%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.694"
$ws.Range("E10").Value = "  +9.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.128"
$ws.Range("E11").Value = "  +30.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.94"
$ws.Range("E12").Value = "  +10.09%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").Value = "4.002.35"
$ws.Range("E14").Value = "  +3.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.77"
$ws.Range("E15").Value = "  +4.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.31"
$ws.Range("E16").Value = "  +4.31%  "

$ws.Range("D17").Value = "3.467.06"
$ws.Range("E17").Value = "  +3.74%  "

$ws.Range("D18").Value = "62.549.25"
$ws.Range("E18").Value = "  +6.05%  "

$ws.Range("E19").Value = "  +0.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000139"
$ws.Range("E21").Value = "  +26.96%  "

$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.16"
$ws.Range("E24").Value = "  +9.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "313.59"
$ws.Range("E25").Value = "  +3.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.22"
$ws.Range("E26").Value = "  -1.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.40"
$ws.Range("E27").Value = "  +5.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.80"
$ws.Range("E29").Value = "  +6.06%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  +7.77%  "

$ws.Range("E31").Value = "  +4.14%  "

$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "44.91"
$ws.Range("E33").Value = "  +11.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.98"
$ws.Range("E34").Value = "  +4.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.67"
$ws.Range("E35").Value = "  +25.14%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0494"
$ws.Range("E37").Value = "  -8.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.74"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("E39").Value = "  +2.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  -6.12%  "

$ws.Range("E42").Value = "  +5.21%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.05"
$ws.Range("E43").Value = "  +7.49%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.126"
$ws.Range("E44").Value = "  +2.85%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "137.57"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.290"
$ws.Range("E46").Value = "  +3.84%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.01"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.27"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.67"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("D50").Value = "2.250.81"
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("D51").Value = "3.807.70"
$ws.Range("E51").Value = "  +4.16%  "

